$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -1.365967021313572
$ws.Range("D2").Value = 0.1857517862905902

$ws.Range("C3").Value = 0.2708532307417135
$ws.Range("D3").Value = 0.7890289110717494

$ws.Range("C4").Value = 0.3809840014954694
$ws.Range("D4").Value = 0.7068693809724627

$ws.Range("C5").Value = -0.7608707682253308
$ws.Range("D5").Value = 0.4548127916673992

$ws.Range("C6").Value = 1.372684434917427
$ws.Range("D6").Value = 0.1836830252947494

$ws.Range("C7").Value = 1.943956698312353
$ws.Range("D7").Value = 0.06479967657682528

$ws.Range("C8").Value = 0.983157800323479
$ws.Range("D8").Value = 0.336219256605387

$ws.Range("C9").Value = 0.1365251797606979
$ws.Range("D9").Value = 0.8926477706250453

$ws.Range("C10").Value = -0.840107012847862
$ws.Range("D10").Value = 0.4098859307663605

$ws.Range("C11").Value = -0.957617674931979
$ws.Range("D11").Value = 0.3486617313706291
